$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3138.102
$ws.Range("J17").Value = 3155.681
$ws.Range("L17").Value = 9467.043
$ws.Range("N17").Value = -9803.043

$ws.Range("H116").Value = 8050.5557
$ws.Range("I116").Value = 8082.0835
$ws.Range("K116").Value = 8082.0835
$ws.Range("M116").Value = -4640.0835

$ws.Range("H131").Value = 4788.9
$ws.Range("J131").Value = 5457
$ws.Range("L131").Value = 16371
$ws.Range("N131").Value = -26451

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 59999
$ws.Range("J23").Value = 59999
$ws.Range("L23").Value = 59999
$ws.Range("N23").Value = -60517

$ws.Range("H61").Value = 3778.3333
$ws.Range("I61").Value = 3261.25
$ws.Range("J61").Value = 4812.5
$ws.Range("K61").Value = 3261.25
$ws.Range("L61").Value = 4812.5
$ws.Range("M61").Value = -3049.25
$ws.Range("N61").Value = -5236.5

$ws.Range("H63").Value = 4172.364
$ws.Range("J63").Value = 4544
$ws.Range("L63").Value = 4544
$ws.Range("N63").Value = -5916

$ws.Range("H66").Value = 4172.364
$ws.Range("J66").Value = 4544
$ws.Range("L66").Value = 22720
$ws.Range("N66").Value = -29584

$ws.Range("H74").Value = 3119
$ws.Range("I74").Value = 2721.2354
$ws.Range("K74").Value = 2721.2354
$ws.Range("M74").Value = -1847.2354

$ws.Range("H77").Value = 3119
$ws.Range("I77").Value = 2721.2354
$ws.Range("K77").Value = 13606.177
$ws.Range("M77").Value = -9238.177

$ws.Range("H97").Value = 3422.3076
$ws.Range("I97").Value = 2030
$ws.Range("K97").Value = 2030
$ws.Range("M97").Value = -1534

$ws.Range("H110").Value = 1611.5
$ws.Range("J110").Value = 2566.6667
$ws.Range("L110").Value = 2566.6667
$ws.Range("N110").Value = -6656.6667

$ws.Range("H132").Value = 3873.8206
$ws.Range("I132").Value = 3719.2856
$ws.Range("J132").Value = 5226
$ws.Range("K132").Value = 11157.8568
$ws.Range("L132").Value = 15678
$ws.Range("M132").Value = -8627.856800000001
$ws.Range("N132").Value = -20738

$ws.Range("H136").Value = 3778.3333
$ws.Range("I136").Value = 3261.25
$ws.Range("J136").Value = 4812.5
$ws.Range("K136").Value = 9783.75
$ws.Range("L136").Value = 14437.5
$ws.Range("M136").Value = -7233.75
$ws.Range("N136").Value = -19537.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7145931
$ws.Range("I134").Value = 8930601
$ws.Range("K134").Value = 26791803
$ws.Range("M134").Value = -26789268

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 34.375
$ws.Range("I7").Value = 32.35
$ws.Range("K7").Value = 32.35
$ws.Range("M7").Value = 80.65000000000001

$ws.Range("H16").Value = 1332
$ws.Range("I16").Value = 1461.3334
$ws.Range("K16").Value = 1461.3334
$ws.Range("M16").Value = -1174.3334

$ws.Range("H31").Value = 4215.1113
$ws.Range("I31").Value = 2032.6666
$ws.Range("J31").Value = 4651.6
$ws.Range("K31").Value = 2032.6666
$ws.Range("L31").Value = 4651.6
$ws.Range("M31").Value = -1737.6666
$ws.Range("N31").Value = -5241.6

$ws.Range("H34").Value = 4215.1113
$ws.Range("I34").Value = 2032.6666
$ws.Range("J34").Value = 4651.6
$ws.Range("K34").Value = 2032.6666
$ws.Range("L34").Value = 4651.6
$ws.Range("M34").Value = -1830.6666
$ws.Range("N34").Value = -5055.6

$ws.Range("H68").Value = 74500
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 74500
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 74500
$ws.Range("N68").Value = -75998
$ws.Range("M68").Value = ""

$ws.Range("H71").Value = 74500
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 74500
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 223500
$ws.Range("N71").Value = -230988
$ws.Range("M71").Value = ""

$ws.Range("H113").Value = 1332
$ws.Range("I113").Value = 1461.3334
$ws.Range("K113").Value = 1461.3334
$ws.Range("M113").Value = 708.6666

$ws.Range("H141").Value = 303866.25
$ws.Range("J141").Value = 390155.34
$ws.Range("L141").Value = 390155.34
$ws.Range("N141").Value = -400515.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2229166.8
$ws.Range("J9").Value = 8571.429
$ws.Range("L9").Value = 25714.287
$ws.Range("N9").Value = -26162.287

$ws.Range("H107").Value = 1110.7894
$ws.Range("I107").Value = 1221.5385
$ws.Range("J107").Value = 870.8333
$ws.Range("K107").Value = 3664.6155
$ws.Range("L107").Value = 2612.4999
$ws.Range("M107").Value = -1744.6155
$ws.Range("N107").Value = -6452.4999

$ws.Range("H129").Value = 1531.6666
$ws.Range("I129").Value = 798.5
$ws.Range("J129").Value = 2998
$ws.Range("K129").Value = 2395.5
$ws.Range("L129").Value = 8994
$ws.Range("M129").Value = 2604.5
$ws.Range("N129").Value = -18994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 74500
$ws.Range("J69").Value = 74500
$ws.Range("L69").Value = 74500
$ws.Range("N69").Value = -75998

$ws.Range("H72").Value = 74500
$ws.Range("J72").Value = 74500
$ws.Range("L72").Value = 223500
$ws.Range("N72").Value = -230988

$ws.Range("H113").Value = 14621.25
$ws.Range("I113").Value = 19929.334
$ws.Range("J113").Value = 11436.4
$ws.Range("K113").Value = 19929.334
$ws.Range("L113").Value = 11436.4
$ws.Range("M113").Value = -17759.334
$ws.Range("N113").Value = -15776.4

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""

$ws.Range("H132").Value = 2797.2942
$ws.Range("I132").Value = 2436.9333
$ws.Range("K132").Value = 7310.7999
$ws.Range("M132").Value = -4780.7999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4299.7144
$ws.Range("I68").Value = 5019.6
$ws.Range("K68").Value = 5019.6
$ws.Range("M68").Value = -4270.6

$ws.Range("H71").Value = 4299.7144
$ws.Range("I71").Value = 5019.6
$ws.Range("K71").Value = 25098
$ws.Range("M71").Value = -21354

$ws.Range("H119").Value = 98989
$ws.Range("J119").Value = 98989
$ws.Range("L119").Value = 98989
$ws.Range("N119").Value = -108665

$ws.Range("H136").Value = 13052.5
$ws.Range("J136").Value = 21285.428
$ws.Range("L136").Value = 63856.284
$ws.Range("N136").Value = -68956.284
